{"js": "const replacements = [\n  [\"2023-03-07 Tuesday\", \"2023-03-08 Wednesday\"],\n  [\"23\u00d773=\", \"14\u00d746=\"],\n  [\"69\u00d785=\", \"66\u00d786=\"],\n  [\"32\u00d778=\", \"81\u00d711=\"],\n  [\"18\u00d742=\", \"20\u00d775=\"],\n  [\"11\u00d730=\", \"61\u00d780=\"],\n  [\"99\u00d720=\", \"96\u00d770=\"],\n  [\"42\u00d725=\", \"20\u00d733=\"],\n  [\"59\u00d726=\", \"48\u00d764=\"],\n  [\"82\u00d778=\", \"39\u00d794=\"],\n  [\"19\u00d797=\", \"87\u00d733=\"],\n  [\"18\u00d775=\", \"84\u00d760=\"],\n  [\"44\u00d739=\", \"57\u00d748=\"],\n  [\"57\u00d718=\", \"70\u00d759=\"],\n  [\"30\u00d772=\", \"80\u00d766=\"],\n  [\"52\u00d721=\", \"19\u00d763=\"],\n  [\"19\u00d767=\", \"21\u00d761=\"],\n  [\"95\u00d717=\", \"48\u00d739=\"],\n  [\"60\u00d799=\", \"78\u00d773=\"],\n  [\"99\u00d793=\", \"10\u00d715=\"],\n  [\"53\u00d746=\", \"56\u00d728=\"],\n  [\"24\u00d722=\", \"57\u00d724=\"],\n  [\"82\u00d712=\", \"45\u00d794=\"],\n  [\"14\u00d731=\", \"84\u00d721=\"],\n  [\"39\u00d751=\", \"35\u00d725=\"],\n  [\"96\u00d724=\", \"20\u00d715=\"],\n  [\"82\u00d775=\", \"91\u00d752=\"],\n  [\"27\u00d728=\", \"34\u00d732=\"],\n  [\"83\u00d794=\", \"98\u00d713=\"],\n  [\"62\u00d749=\", \"90\u00d733=\"],\n  [\"61\u00d734=\", \"87\u00d738=\"],\n  [\"70\u00d716=\", \"47\u00d760=\"],\n  [\"97\u00d728=\", \"76\u00d796=\"],\n  [\"84\u00d745=\", \"62\u00d715=\"],\n  [\"75\u00d776=\", \"77\u00d763=\"],\n  [\"72\u00d741=\", \"96\u00d786=\"],\n  [\"50\u00d792=\", \"45\u00d734=\"],\n  [\"77\u00d774=\", \"69\u00d753=\"],\n  [\"20\u00d740=\", \"48\u00d752=\"],\n  [\"73\u00d712=\", \"95\u00d785=\"],\n  [\"71\u00d789=\", \"27\u00d735=\"],\n  [\"71\u00d793=\", \"26\u00d747=\"],\n  [\"44\u00d734=\", \"81\u00d789=\"],\n  [\"22\u00d714=\", \"95\u00d714=\"],\n  [\"54\u00d717=\", \"86\u00d779=\"],\n  [\"75\u00d728=\", \"47\u00d785=\"],\n  [\"72\u00d764=\", \"91\u00d724=\"],\n  [\"25\u00d755=\", \"62\u00d724=\"],\n  [\"72\u00d767=\", \"58\u00d740=\"],\n  [\"12\u00d749=\", \"55\u00d710=\"],\n  [\"60\u00d713=\", \"45\u00d790=\"],\n  [\"26\u00d777=\", \"32\u00d735=\"],\n  [\"100\u00d782=\", \"57\u00d716=\"],\n  [\"12\u00d720=\", \"89\u00d719=\"],\n  [\"81\u00d744=\", \"41\u00d752=\"],\n  [\"87\u00d774=\", \"12\u00d713=\"],\n  [\"39\u00d792=\", \"61\u00d750=\"],\n  [\"62\u00d774=\", \"95\u00d733=\"],\n  [\"21\u00d7100=\", \"53\u00d775=\"],\n  [\"38\u00d744=\", \"47\u00d723=\"],\n  [\"47\u00d774=\", \"10\u00d736=\"],\n  [\"69\u00d755=\", \"15\u00d772=\"],\n  [\"50\u00d740=\", \"37\u00d726=\"],\n  [\"43\u00d785=\", \"87\u00d743=\"],\n  [\"92\u00d724=\", \"69\u00d720=\"],\n  [\"37\u00d784=\", \"79\u00d784=\"],\n  [\"55\u00d750=\", \"28\u00d763=\"],\n  [\"62\u00d713=\", \"72\u00d780=\"],\n  [\"45\u00d716=\", \"86\u00d784=\"],\n  [\"57\u00d713=\", \"60\u00d767=\"],\n  [\"36\u00d786=\", \"17\u00d721=\"],\n  [\"65\u00d748=\", \"46\u00d793=\"],\n  [\"70\u00d751=\", \"36\u00d765=\"],\n  [\"65\u00d743=\", \"100\u00d786=\"],\n  [\"20\u00d752=\", \"52\u00d752=\"],\n  [\"73\u00d780=\", \"68\u00d740=\"],\n  [\"19\u00d759=\", \"80\u00d738=\"],\n  [\"45\u00d773=\", \"30\u00d728=\"],\n  [\"22\u00d738=\", \"63\u00d785=\"],\n  [\"14\u00d771=\", \"87\u00d731=\"],\n  [\"63\u00d746=\", \"31\u00d710=\"],\n  [\"99\u00d782=\", \"85\u00d752=\"],\n  [\"97\u00d734=\", \"88\u00d737=\"],\n  [\"12\u00d736=\", \"28\u00d720=\"],\n  [\"53\u00d795=\", \"32\u00d748=\"],\n  [\"70\u00d794=\", \"56\u00d754=\"],\n  [\"22\u00d719=\", \"29\u00d725=\"],\n  [\"55\u00d795=\", \"17\u00d770=\"],\n  [\"97\u00d745=\", \"25\u00d740=\"],\n  [\"88\u00d789=\", \"65\u00d779=\"],\n  [\"32\u00d771=\", \"50\u00d733=\"],\n  [\"31\u00d790=\", \"67\u00d748=\"],\n  [\"39\u00d759=\", \"39\u00d728=\"],\n  [\"96\u00d761=\", \"62\u00d771=\"],\n  [\"69\u00d748=\", \"37\u00d785=\"],\n  [\"83\u00d781=\", \"71\u00d713=\"],\n  [\"11\u00d761=\", \"80\u00d797=\"],\n  [\"70\u00d731=\", \"95\u00d713=\"],\n  [\"23\u00d777=\", \"33\u00d767=\"],\n  [\"21\u00d759=\", \"80\u00d775=\"],\n  [\"33\u00d785=\", \"24\u00d734=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2023-03-07 Tuesday\", \"2023-03-08 Wednesday\"),\n  @(\"23\u00d773=\", \"14\u00d746=\"),\n  @(\"69\u00d785=\", \"66\u00d786=\"),\n  @(\"32\u00d778=\", \"81\u00d711=\"),\n  @(\"18\u00d742=\", \"20\u00d775=\"),\n  @(\"11\u00d730=\", \"61\u00d780=\"),\n  @(\"99\u00d720=\", \"96\u00d770=\"),\n  @(\"42\u00d725=\", \"20\u00d733=\"),\n  @(\"59\u00d726=\", \"48\u00d764=\"),\n  @(\"82\u00d778=\", \"39\u00d794=\"),\n  @(\"19\u00d797=\", \"87\u00d733=\"),\n  @(\"18\u00d775=\", \"84\u00d760=\"),\n  @(\"44\u00d739=\", \"57\u00d748=\"),\n  @(\"57\u00d718=\", \"70\u00d759=\"),\n  @(\"30\u00d772=\", \"80\u00d766=\"),\n  @(\"52\u00d721=\", \"19\u00d763=\"),\n  @(\"19\u00d767=\", \"21\u00d761=\"),\n  @(\"95\u00d717=\", \"48\u00d739=\"),\n  @(\"60\u00d799=\", \"78\u00d773=\"),\n  @(\"99\u00d793=\", \"10\u00d715=\"),\n  @(\"53\u00d746=\", \"56\u00d728=\"),\n  @(\"24\u00d722=\", \"57\u00d724=\"),\n  @(\"82\u00d712=\", \"45\u00d794=\"),\n  @(\"14\u00d731=\", \"84\u00d721=\"),\n  @(\"39\u00d751=\", \"35\u00d725=\"),\n  @(\"96\u00d724=\", \"20\u00d715=\"),\n  @(\"82\u00d775=\", \"91\u00d752=\"),\n  @(\"27\u00d728=\", \"34\u00d732=\"),\n  @(\"83\u00d794=\", \"98\u00d713=\"),\n  @(\"62\u00d749=\", \"90\u00d733=\"),\n  @(\"61\u00d734=\", \"87\u00d738=\"),\n  @(\"70\u00d716=\", \"47\u00d760=\"),\n  @(\"97\u00d728=\", \"76\u00d796=\"),\n  @(\"84\u00d745=\", \"62\u00d715=\"),\n  @(\"75\u00d776=\", \"77\u00d763=\"),\n  @(\"72\u00d741=\", \"96\u00d786=\"),\n  @(\"50\u00d792=\", \"45\u00d734=\"),\n  @(\"77\u00d774=\", \"69\u00d753=\"),\n  @(\"20\u00d740=\", \"48\u00d752=\"),\n  @(\"73\u00d712=\", \"95\u00d785=\"),\n  @(\"71\u00d789=\", \"27\u00d735=\"),\n  @(\"71\u00d793=\", \"26\u00d747=\"),\n  @(\"44\u00d734=\", \"81\u00d789=\"),\n  @(\"22\u00d714=\", \"95\u00d714=\"),\n  @(\"54\u00d717=\", \"86\u00d779=\"),\n  @(\"75\u00d728=\", \"47\u00d785=\"),\n  @(\"72\u00d764=\", \"91\u00d724=\"),\n  @(\"25\u00d755=\", \"62\u00d724=\"),\n  @(\"72\u00d767=\", \"58\u00d740=\"),\n  @(\"12\u00d749=\", \"55\u00d710=\"),\n  @(\"60\u00d713=\", \"45\u00d790=\"),\n  @(\"26\u00d777=\", \"32\u00d735=\"),\n  @(\"100\u00d782=\", \"57\u00d716=\"),\n  @(\"12\u00d720=\", \"89\u00d719=\"),\n  @(\"81\u00d744=\", \"41\u00d752=\"),\n  @(\"87\u00d774=\", \"12\u00d713=\"),\n  @(\"39\u00d792=\", \"61\u00d750=\"),\n  @(\"62\u00d774=\", \"95\u00d733=\"),\n  @(\"21\u00d7100=\", \"53\u00d775=\"),\n  @(\"38\u00d744=\", \"47\u00d723=\"),\n  @(\"47\u00d774=\", \"10\u00d736=\"),\n  @(\"69\u00d755=\", \"15\u00d772=\"),\n  @(\"50\u00d740=\", \"37\u00d726=\"),\n  @(\"43\u00d785=\", \"87\u00d743=\"),\n  @(\"92\u00d724=\", \"69\u00d720=\"),\n  @(\"37\u00d784=\", \"79\u00d784=\"),\n  @(\"55\u00d750=\", \"28\u00d763=\"),\n  @(\"62\u00d713=\", \"72\u00d780=\"),\n  @(\"45\u00d716=\", \"86\u00d784=\"),\n  @(\"57\u00d713=\", \"60\u00d767=\"),\n  @(\"36\u00d786=\", \"17\u00d721=\"),\n  @(\"65\u00d748=\", \"46\u00d793=\"),\n  @(\"70\u00d751=\", \"36\u00d765=\"),\n  @(\"65\u00d743=\", \"100\u00d786=\"),\n  @(\"20\u00d752=\", \"52\u00d752=\"),\n  @(\"73\u00d780=\", \"68\u00d740=\"),\n  @(\"19\u00d759=\", \"80\u00d738=\"),\n  @(\"45\u00d773=\", \"30\u00d728=\"),\n  @(\"22\u00d738=\", \"63\u00d785=\"),\n  @(\"14\u00d771=\", \"87\u00d731=\"),\n  @(\"63\u00d746=\", \"31\u00d710=\"),\n  @(\"99\u00d782=\", \"85\u00d752=\"),\n  @(\"97\u00d734=\", \"88\u00d737=\"),\n  @(\"12\u00d736=\", \"28\u00d720=\"),\n  @(\"53\u00d795=\", \"32\u00d748=\"),\n  @(\"70\u00d794=\", \"56\u00d754=\"),\n  @(\"22\u00d719=\", \"29\u00d725=\"),\n  @(\"55\u00d795=\", \"17\u00d770=\"),\n  @(\"97\u00d745=\", \"25\u00d740=\"),\n  @(\"88\u00d789=\", \"65\u00d779=\"),\n  @(\"32\u00d771=\", \"50\u00d733=\"),\n  @(\"31\u00d790=\", \"67\u00d748=\"),\n  @(\"39\u00d759=\", \"39\u00d728=\"),\n  @(\"96\u00d761=\", \"62\u00d771=\"),\n  @(\"69\u00d748=\", \"37\u00d785=\"),\n  @(\"83\u00d781=\", \"71\u00d713=\"),\n  @(\"11\u00d761=\", \"80\u00d797=\"),\n  @(\"70\u00d731=\", \"95\u00d713=\"),\n  @(\"23\u00d777=\", \"33\u00d767=\"),\n  @(\"21\u00d759=\", \"80\u00d775=\"),\n  @(\"33\u00d785=\", \"24\u00d734=\"),\n)\n\nforeach ($pair in $replacements) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $range = $d.Content\n  $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}"}
